$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks Description")

# --- Complete the "site.md" Description block (rows 44-46), mirroring the
#     pattern used by the other sections (copy style+value from the Theme
#     section's matching Description block). ---
$ws.Range("F12").Copy($ws.Range("F44"))
$ws.Range("G23").Copy($ws.Range("G45"))
$ws.Range("G24").Copy($ws.Range("G46"))

# --- Split the old row 50 (E50 "Publishing" header + F50 "Resources:")
#     into its own header row (47) and Resources row (48). ---
$ws.Range("E5:F5").Copy($ws.Range("E47:F47"))
$ws.Range("E47").Value = "Publishing"
$ws.Range("F6").Copy($ws.Range("F48"))

# --- Move the old row 51 (Staging resource link) down to row 49. ---
$ws.Range("G51:I51").Copy($ws.Range("G49:I49"))

# --- New Description: block for Publishing (rows 50-53). ---
$ws.Range("F12").Copy($ws.Range("F50"))
$ws.Range("G51").Value = "Create a production site"
$ws.Range("I51").ClearContents()
$ws.Range("G52").Value = "Create a staging site"
$ws.Range("G53").Value = "One person in charge of"

# --- Move old row 53 ("Build Site" header) down to row 54. ---
$ws.Range("E54").Value = "Build Site"
$ws.Range("E50").ClearContents()
$ws.Range("E53").ClearContents()

# --- View state: "Tasks Description" becomes the active/selected sheet,
#     with E53 selected; "Gattn Chart" is no longer the selected tab. ---
$ws.Activate()
$ws.Range("E53").Select()

$ws2 = $wb.Worksheets.Item("Gattn Chart")
$ws2.Range("F18").Select()
$ws.Activate()
